$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.436.34"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.574.90"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.28"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3767"
$ws.Range("E7").Value = "  +2.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.93"
$ws.Range("E8").Value = "  +0.66%  "

# Row 9
$ws.Range("E9").Value = "  +1.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.165"
$ws.Range("E10").Value = "  -0.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07675"
$ws.Range("E11").Value = "  +1.57%  "

# Row 12
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.29"
$ws.Range("E13").Value = "  +1.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.003"
$ws.Range("E14").Value = "  -0.62%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.935"
$ws.Range("E15").Value = "  +1.39%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.575.41"
$ws.Range("E16").Value = "  +0.07%  "

# Row 17
$ws.Range("E17").Value = "  -0.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.42"
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06771"
$ws.Range("E19").Value = "  +1.10%  "

# Row 20
$ws.Range("E20").Value = "  +0.13%  "

# Row 21
$ws.Range("E21").Value = "  +3.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.245"

# Row 23
$ws.Range("E23").Value = "  +1.13%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.426"
$ws.Range("E24").Value = "  +1.10%  "

# Row 25
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.437.73"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.750"
$ws.Range("E26").Value = "  -6.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.37"
$ws.Range("E27").Value = "  +2.72%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "145.89"
$ws.Range("E28").Value = "  -0.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.034"
$ws.Range("E29").Value = "  +1.77%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.28"
$ws.Range("E30").Value = "  +1.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.755.41"
$ws.Range("E31").Value = "  +0.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.241"
$ws.Range("E32").Value = "  -0.31%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.014"
$ws.Range("E33").Value = "  +3.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.016"
$ws.Range("E34").Value = "  +1.96%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.04"
$ws.Range("E35").Value = "  -3.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08578"
$ws.Range("E36").Value = "  +1.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02562"
$ws.Range("E37").Value = "  +1.37%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2316"
$ws.Range("E38").Value = "  +0.81%  "

# Row 39
$ws.Range("E39").Value = "  +1.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.343"
$ws.Range("E40").Value = "  +8.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.471"
$ws.Range("E41").Value = "  -0.91%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.61"
$ws.Range("E42").Value = "  -1.35%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6469"
$ws.Range("E43").Value = "  +1.47%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.11"
$ws.Range("E44").Value = "  -2.84%  "

# Row 45
$ws.Range("E45").Value = "  +0.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6037"
$ws.Range("E46").Value = "  +0.44%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.801"
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.302"
$ws.Range("E48").Value = "  +9.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.095"
$ws.Range("E49").Value = "  -0.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.52"
$ws.Range("E50").Value = "  +3.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07335"
$ws.Range("E51").Value = "  +0.95%  "
